# Applies the 2024-11-22 05:26:25 crypto live-data refresh described in the commit diff.
# Updates price/market-cap/volume/change figures on all three sheets, plus the
# Pepe/Polkadot (rows 22-23) and Filecoin/Mantle (rows 50-51) re-ranking swaps.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Top 50 Cryptocurrencies ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 3).Value = 98817
$ws1.Cells.Item(2, 4).Value = 1953772577503
$ws1.Cells.Item(2, 5).Value = 111602031886
$ws1.Cells.Item(2, 6).Value = 1.18963

$ws1.Cells.Item(3, 3).Value = 3383.02
$ws1.Cells.Item(3, 4).Value = 407043413496
$ws1.Cells.Item(3, 5).Value = 56933584709
$ws1.Cells.Item(3, 6).Value = 7.99777

$ws1.Cells.Item(4, 4).Value = 130788831244
$ws1.Cells.Item(4, 5).Value = 95661940820
$ws1.Cells.Item(4, 6).Value = -0.01272

$ws1.Cells.Item(5, 3).Value = 261.91
$ws1.Cells.Item(5, 4).Value = 124075895781
$ws1.Cells.Item(5, 5).Value = 14990065171
$ws1.Cells.Item(5, 6).Value = 8.63114

$ws1.Cells.Item(6, 3).Value = 635.51
$ws1.Cells.Item(6, 4).Value = 92676534988
$ws1.Cells.Item(6, 5).Value = 2458154483
$ws1.Cells.Item(6, 6).Value = 3.874

$ws1.Cells.Item(7, 3).Value = 1.39
$ws1.Cells.Item(7, 4).Value = 78917385028
$ws1.Cells.Item(7, 5).Value = 17890002015
$ws1.Cells.Item(7, 6).Value = 24.17602

$ws1.Cells.Item(8, 3).Value = 0.395961
$ws1.Cells.Item(8, 4).Value = 58130923980
$ws1.Cells.Item(8, 5).Value = 9854138290
$ws1.Cells.Item(8, 6).Value = 1.94918

$ws1.Cells.Item(9, 4).Value = 38297868428
$ws1.Cells.Item(9, 5).Value = 11136208756
$ws1.Cells.Item(9, 6).Value = -0.00488

$ws1.Cells.Item(10, 3).Value = 3383.35
$ws1.Cells.Item(10, 4).Value = 33126025781
$ws1.Cells.Item(10, 5).Value = 144753180
$ws1.Cells.Item(10, 6).Value = 8.165369999999999

$ws1.Cells.Item(11, 3).Value = 0.882933
$ws1.Cells.Item(11, 4).Value = 31575973787
$ws1.Cells.Item(11, 5).Value = 3537099617
$ws1.Cells.Item(11, 6).Value = 10.89773

$ws1.Cells.Item(12, 3).Value = 0.200395
$ws1.Cells.Item(12, 4).Value = 17289684171
$ws1.Cells.Item(12, 5).Value = 1081990781
$ws1.Cells.Item(12, 6).Value = 1.38293

$ws1.Cells.Item(13, 3).Value = 36.4
$ws1.Cells.Item(13, 4).Value = 14881406931
$ws1.Cells.Item(13, 5).Value = 1037701522
$ws1.Cells.Item(13, 6).Value = 6.53633

$ws1.Cells.Item(14, 3).Value = 0.00002499
$ws1.Cells.Item(14, 4).Value = 14705498503
$ws1.Cells.Item(14, 5).Value = 1610868051
$ws1.Cells.Item(14, 6).Value = 2.95005

$ws1.Cells.Item(15, 3).Value = 3998.89
$ws1.Cells.Item(15, 4).Value = 14427358496
$ws1.Cells.Item(15, 5).Value = 93745940
$ws1.Cells.Item(15, 6).Value = 7.34263

$ws1.Cells.Item(16, 3).Value = 98913
$ws1.Cells.Item(16, 4).Value = 14395882737
$ws1.Cells.Item(16, 5).Value = 842373855
$ws1.Cells.Item(16, 6).Value = 1.97696

$ws1.Cells.Item(17, 3).Value = 5.56
$ws1.Cells.Item(17, 4).Value = 14140396529
$ws1.Cells.Item(17, 5).Value = 637922301
$ws1.Cells.Item(17, 6).Value = 3.21844

$ws1.Cells.Item(18, 3).Value = 3.61
$ws1.Cells.Item(18, 4).Value = 10263931784
$ws1.Cells.Item(18, 5).Value = 1893045491
$ws1.Cells.Item(18, 6).Value = 0.7609

$ws1.Cells.Item(19, 3).Value = 496.39
$ws1.Cells.Item(19, 4).Value = 9820858921
$ws1.Cells.Item(19, 5).Value = 1972327901
$ws1.Cells.Item(19, 6).Value = -3.06361

$ws1.Cells.Item(20, 3).Value = 3383.89
$ws1.Cells.Item(20, 4).Value = 9627349785
$ws1.Cells.Item(20, 5).Value = 1238335912
$ws1.Cells.Item(20, 6).Value = 8.11872

$ws1.Cells.Item(21, 3).Value = 15.28
$ws1.Cells.Item(21, 4).Value = 9558865243
$ws1.Cells.Item(21, 5).Value = 1247144986
$ws1.Cells.Item(21, 6).Value = 4.32433

$ws1.Cells.Item(22, 1).Value = 'Polkadot'
$ws1.Cells.Item(22, 2).Value = 'dot'
$ws1.Cells.Item(22, 3).Value = 6.23
$ws1.Cells.Item(22, 4).Value = 8966606443
$ws1.Cells.Item(22, 5).Value = 823772180
$ws1.Cells.Item(22, 6).Value = 8.80635

$ws1.Cells.Item(23, 1).Value = 'Pepe'
$ws1.Cells.Item(23, 2).Value = 'pepe'
$ws1.Cells.Item(23, 3).Value = 0.00002122
$ws1.Cells.Item(23, 4).Value = 8928147587
$ws1.Cells.Item(23, 5).Value = 6810977654
$ws1.Cells.Item(23, 6).Value = 8.61903

$ws1.Cells.Item(24, 3).Value = 0.283327
$ws1.Cells.Item(24, 4).Value = 8497950433
$ws1.Cells.Item(24, 5).Value = 2292477997
$ws1.Cells.Item(24, 6).Value = 17.93069

$ws1.Cells.Item(25, 3).Value = 8.76
$ws1.Cells.Item(25, 4).Value = 8102501993
$ws1.Cells.Item(25, 5).Value = 3436944
$ws1.Cells.Item(25, 6).Value = 2.93876

$ws1.Cells.Item(26, 3).Value = 5.81
$ws1.Cells.Item(26, 4).Value = 7068507828
$ws1.Cells.Item(26, 5).Value = 1011665226
$ws1.Cells.Item(26, 6).Value = 4.38588

$ws1.Cells.Item(27, 3).Value = 90.5
$ws1.Cells.Item(27, 4).Value = 6802345942
$ws1.Cells.Item(27, 5).Value = 1416854551
$ws1.Cells.Item(27, 6).Value = 4.09299

$ws1.Cells.Item(28, 3).Value = 12.14
$ws1.Cells.Item(28, 4).Value = 6466066532
$ws1.Cells.Item(28, 5).Value = 865453610
$ws1.Cells.Item(28, 6).Value = 4.19963

$ws1.Cells.Item(29, 3).Value = 3562.66
$ws1.Cells.Item(29, 4).Value = 6191810117
$ws1.Cells.Item(29, 5).Value = 102355778
$ws1.Cells.Item(29, 6).Value = 8.005269999999999

$ws1.Cells.Item(30, 4).Value = 5647368459
$ws1.Cells.Item(30, 5).Value = 858974039
$ws1.Cells.Item(30, 6).Value = 6.04675

$ws1.Cells.Item(31, 3).Value = 0.199016
$ws1.Cells.Item(31, 4).Value = 5400717796
$ws1.Cells.Item(31, 5).Value = 124320291
$ws1.Cells.Item(31, 6).Value = 13.2041

$ws1.Cells.Item(32, 3).Value = 0.998515
$ws1.Cells.Item(32, 4).Value = 5222825181
$ws1.Cells.Item(32, 5).Value = 16300713
$ws1.Cells.Item(32, 6).Value = -0.78412

$ws1.Cells.Item(33, 3).Value = 0.134483
$ws1.Cells.Item(33, 4).Value = 5137474725
$ws1.Cells.Item(33, 5).Value = 874776333
$ws1.Cells.Item(33, 6).Value = 6.51209

$ws1.Cells.Item(34, 3).Value = 9.69
$ws1.Cells.Item(34, 4).Value = 4598987074
$ws1.Cells.Item(34, 5).Value = 273681702
$ws1.Cells.Item(34, 6).Value = 6.55566

$ws1.Cells.Item(35, 3).Value = 27.94
$ws1.Cells.Item(35, 4).Value = 4178731627
$ws1.Cells.Item(35, 5).Value = 884770184
$ws1.Cells.Item(35, 6).Value = 5.19909

$ws1.Cells.Item(36, 3).Value = 0.00005199
$ws1.Cells.Item(36, 4).Value = 3901948320
$ws1.Cells.Item(36, 5).Value = 1688917428
$ws1.Cells.Item(36, 6).Value = 1.89678

$ws1.Cells.Item(37, 3).Value = 7.38
$ws1.Cells.Item(37, 4).Value = 3815241912
$ws1.Cells.Item(37, 5).Value = 434605828
$ws1.Cells.Item(37, 6).Value = -0.28555

$ws1.Cells.Item(38, 3).Value = 0.151407
$ws1.Cells.Item(38, 4).Value = 3808608303
$ws1.Cells.Item(38, 5).Value = 152035278
$ws1.Cells.Item(38, 6).Value = -0.54395

$ws1.Cells.Item(39, 3).Value = 0.47192
$ws1.Cells.Item(39, 4).Value = 3760540435
$ws1.Cells.Item(39, 5).Value = 489495021
$ws1.Cells.Item(39, 6).Value = 7.29274

$ws1.Cells.Item(40, 3).Value = 508.17
$ws1.Cells.Item(40, 4).Value = 3750000473
$ws1.Cells.Item(40, 5).Value = 286015428
$ws1.Cells.Item(40, 6).Value = 3.47623

$ws1.Cells.Item(41, 4).Value = 3687728228
$ws1.Cells.Item(41, 5).Value = 223922865
$ws1.Cells.Item(41, 6).Value = 0.06525

$ws1.Cells.Item(42, 4).Value = 3582963813
$ws1.Cells.Item(42, 5).Value = 33496946
$ws1.Cells.Item(42, 6).Value = 2.63912

$ws1.Cells.Item(43, 4).Value = 3439078996
$ws1.Cells.Item(43, 5).Value = 154364886
$ws1.Cells.Item(43, 6).Value = -0.00669

$ws1.Cells.Item(44, 3).Value = 3.4
$ws1.Cells.Item(44, 4).Value = 3396851368
$ws1.Cells.Item(44, 5).Value = 1281631615
$ws1.Cells.Item(44, 6).Value = 6.78171

$ws1.Cells.Item(45, 3).Value = 3.72
$ws1.Cells.Item(45, 4).Value = 3361552390
$ws1.Cells.Item(45, 5).Value = 301693729
$ws1.Cells.Item(45, 6).Value = 2.64148

$ws1.Cells.Item(46, 4).Value = 3337822225
$ws1.Cells.Item(46, 5).Value = 488298136
$ws1.Cells.Item(46, 6).Value = 2.42598

$ws1.Cells.Item(47, 3).Value = 0.787967
$ws1.Cells.Item(47, 4).Value = 3226906040
$ws1.Cells.Item(47, 5).Value = 1668181094
$ws1.Cells.Item(47, 6).Value = 12.37254

$ws1.Cells.Item(48, 3).Value = 160.67
$ws1.Cells.Item(48, 4).Value = 2963867990
$ws1.Cells.Item(48, 5).Value = 83352737
$ws1.Cells.Item(48, 6).Value = -1.78004

$ws1.Cells.Item(49, 3).Value = 1.95
$ws1.Cells.Item(49, 4).Value = 2933339796
$ws1.Cells.Item(49, 5).Value = 368223956
$ws1.Cells.Item(49, 6).Value = 0.21023

$ws1.Cells.Item(50, 1).Value = 'Mantle'
$ws1.Cells.Item(50, 2).Value = 'mnt'
$ws1.Cells.Item(50, 3).Value = 0.846624
$ws1.Cells.Item(50, 4).Value = 2850586086
$ws1.Cells.Item(50, 5).Value = 183686556
$ws1.Cells.Item(50, 6).Value = 15.21575

$ws1.Cells.Item(51, 1).Value = 'Filecoin'
$ws1.Cells.Item(51, 2).Value = 'fil'
$ws1.Cells.Item(51, 3).Value = 4.71
$ws1.Cells.Item(51, 4).Value = 2828845966
$ws1.Cells.Item(51, 5).Value = 583490036
$ws1.Cells.Item(51, 6).Value = 7.03179

# --- Sheet 2: Top 5 by Market Cap ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 2).Value = 1953772577503
$ws2.Cells.Item(3, 2).Value = 407043413496
$ws2.Cells.Item(4, 2).Value = 130788831244
$ws2.Cells.Item(5, 2).Value = 124075895781
$ws2.Cells.Item(6, 2).Value = 92676534988

# --- Sheet 3: Summary ---
# Leading "'" forces text entry so Excel does not auto-coerce the
# "$"-prefixed average-price string into a currency number.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 2).Value = "'" + '$4355.88'
$ws3.Cells.Item(3, 2).Value = "'" + 'XRP (24.18%)'
$ws3.Cells.Item(4, 2).Value = "'" + 'Bitcoin Cash (-3.06%)'
